# TC05_C3DC_phs000463_DiagnClasfSystem-ICD-O-3.xlsx
# "Updated queries for C3DC first half testcases."
#
# The sample SQL queries stored in column B/C of Sheet1 joined df_participant /
# df_diagnoses / df_treatments / df_treatment_resp / df_survival / df_reference_files
# using the generic "id" column. The join columns were renamed to their
# fully-qualified equivalents (study_id / participant_id) everywhere the join
# appears across every query cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Sql($text) {
    $t = $text
    $t = $t -replace [regex]::Escape('df_participant prt ON std.id = prt."study.id"'), 'df_participant prt ON std.study_id = prt."study.study_id"'
    $t = $t -replace [regex]::Escape('df_diagnoses dgn ON prt.id = dgn."participant.id"'), 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"'
    $t = $t -replace [regex]::Escape('df_treatments trt ON prt.id = trt."participant.id"'), 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"'
    $t = $t -replace [regex]::Escape('df_treatment_resp trr ON prt.id = trr."participant.id"'), 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"'
    $t = $t -replace [regex]::Escape('df_survival srv ON prt.id = srv."participant.id"'), 'df_survival srv ON prt.participant_id = srv."participant.participant_id"'
    $t = $t -replace [regex]::Escape('df_reference_files rfs ON std.id = rfs."study.id"'), 'df_reference_files rfs ON std.study_id = rfs."study.study_id"'
    return $t
}

# Every cell on the sheet that holds one of the sample SQL queries
# (StudiesTab/ParticipantsTab/DiagnosisTab/TreatmentTab/TreatmentRespTab/
# SurvivalTab query cells, plus the aggregate StatQuery cell).
$queryCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($cell in $queryCells) {
    $orig = $ws.Range($cell).Value2
    $ws.Range($cell).Value = Fix-Sql $orig
}

# The workbook was left with the cursor scrolled to row 7 with C7 selected.
$ws.Range("C7").Select()
